# The "Elitism selection, agents = 2" test run (row 38, A38 = Test ID 36)
# was a duplicate/erroneous entry. Remove it entirely so every later test
# run shifts up by one row, matching the corrected test template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(38).Delete()

# Leave the same selection the author ended up with after tidying the
# template (the remaining block of elitism/tournament/etc. test rows).
$ws.Range("A37:A62").Select()
